$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 151, pushing the existing
# rows 151-168 down to 152-169 (matching dimension growing to A1:T169).
$ws.Range("A151").EntireRow.Insert()

# Populate the newly inserted row 151 with the new weekly price record.
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C151").Value = "Los Lagos"
$ws.Range("D151").Value = 45132
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100104
$ws.Range("H151").Value = "Frutos de pepita"
$ws.Range("I151").Value = 100104003
$ws.Range("J151").Value = "Membrillo"
$ws.Range("K151").Value = "Champion"
$ws.Range("L151").Value = "Primera"
$ws.Range("M151").Value = 100
$ws.Range("N151").Value = 14000
$ws.Range("O151").Value = 14000
$ws.Range("P151").Value = 14000
$ws.Range("Q151").Value = "$/caja 18 kilos empedrada"
$ws.Range("R151").Value = "Región de O'Higgins"
$ws.Range("S151").Value = 778
$ws.Range("T151").Value = 18
